$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "37.844.21"
$ws.Range("E2").Value = "  -0.11%  "
Set-TextCell "D3" "2.084.51"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell "D5" "233.63"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  -0.12%  "
Set-TextCell "D7" "59.13"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("E8").Value = "  -0.06%  "
Set-TextCell "D9" "0.396"
$ws.Range("E9").Value = "  +1.84%  "
Set-TextCell "D10" "0.0788"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +2.73%  "
Set-TextCell "D12" "2.393.40"
$ws.Range("E12").Value = "  +0.17%  "
Set-TextCell "D13" "14.77"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +1.21%  "
Set-TextCell "D16" "5.34"
$ws.Range("E16").Value = "  +1.53%  "
Set-TextCell "D17" "2.049.90"
$ws.Range("E17").Value = "  -2.03%  "
Set-TextCell "D18" "37.777.10"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  +0.17%  "
Set-TextCell "D20" "71.67"
$ws.Range("E20").Value = "  +1.18%  "
Set-TextCell "D21" "0.0₃0849"
$ws.Range("E21").Value = "  +3.07%  "
Set-TextCell "D22" "228.14"
$ws.Range("E22").Value = "  -0.19%  "
Set-TextCell "D23" "0.999"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D26" "9.61"
$ws.Range("E26").Value = "  +7.30%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D27" "171.14"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -2.09%  "
Set-TextCell "D29" "1.42"
$ws.Range("E29").Value = "  -1.08%  "
Set-TextCell "D30" "19.52"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +1.90%  "
Set-TextCell "D32" "4.73"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  -0.39%  "
Set-TextCell "D36" "3.43"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  -1.47%  "
Set-TextCell "D41" "99.34"
$ws.Range("E41").Value = "  +2.05%  "
Set-TextCell "D42" "17.21"
$ws.Range("E42").Value = "  +9.44%  "
$ws.Range("E43").Value = "  +2.41%  "
Set-TextCell "D44" "2.90"
$ws.Range("E44").Value = "  -1.15%  "
Set-TextCell "D45" "1.450.09"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("E48").Value = "  +0.69%  "
Set-TextCell "D49" "7.39"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -1.00%  "
Set-TextCell "D51" "2.277.59"
$ws.Range("E51").Value = "  -0.32%  "
